# Adds "alpha" (E) and "pmax" (F) columns to each data sheet, with
# per-sheet constant values for rows 2-25, formatted as Scientific (0.00E+00)
# for the alpha column (except on the Crypto sheet, which keeps General).
# Also updates each sheet's selection to mirror the authored workbook state.

$wb = $excel.ActiveWorkbook

$sheetsInfo = @(
    @{ Index = 1; Alpha = 0.0000487589400000000011020130641359315859517664648592472076416015625; Pmax = 0.1216388999999999942502881822292692959308624267578125; StyleAlpha = $true;  ActiveCell = "H6";  Sqref = "H6" },
    @{ Index = 2; Alpha = 0.000084925829999999994583508022838458373371395282447338104248046875;   Pmax = 0.2367559999999999942321693424673867411911487579345703125; StyleAlpha = $true;  ActiveCell = "I9";  Sqref = "I9" },
    @{ Index = 3; Alpha = 0.00009941624000000000142308886896813646671944297850131988525390625;    Pmax = 0.08416310999999999931109329054379486478865146636962890625; StyleAlpha = $true;  ActiveCell = "I14"; Sqref = "I14" },
    @{ Index = 4; Alpha = 0.0000496493099999999998295750758270372671177028678357601165771484375;  Pmax = 0.2296909000000000034003022619799594394862651824951171875; StyleAlpha = $true;  ActiveCell = "E2";  Sqref = "E2:F25" },
    @{ Index = 5; Alpha = 0.00010538459999999999386170734805290294389124028384685516357421875;    Pmax = 0.1250165000000000026236790517941699363291263580322265625; StyleAlpha = $false; ActiveCell = "H11"; Sqref = "H11" }
)

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Index)

    # Header row
    $ws.Range("E1").Value = "alpha"
    $ws.Range("F1").Value = "pmax"

    # Data rows 2-25: constant alpha/pmax values per sheet
    $ws.Range("E2:E25").Value = $info.Alpha
    $ws.Range("F2:F25").Value = $info.Pmax

    if ($info.StyleAlpha) {
        $ws.Range("E2:E25").NumberFormat = "0.00E+00"
    }

    # Update the saved selection to match the authored state
    [void]$ws.Range($info.ActiveCell).Select()
    [void]$ws.Range($info.Sqref).Select()
}
